# Generate Report for Handoff
#
# Updates the localization-status report after a handoff run completed for the
# "Ready for handoff" documents (rows 4-7 of the zh-cn and de-de sheets):
#   - Priority moves from "low" to "ht" (handed-off / high-priority triage)
#   - The Latest Handoff Datetime is refreshed to the new handoff timestamp

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Priority column (E) for the four "Ready for handoff" rows switches from low -> ht
$zhcn.Range("E4:E7").Value = "ht"
$dede.Range("E4:E7").Value = "ht"

# Latest Handoff Datetime column (H) for those same rows is refreshed
$zhcn.Range("H4:H7").Value = "2016-09-07 10:44:22"
$dede.Range("H4:H7").Value = "2016-09-07 10:44:28"

# The Overview sheet mirrors the de-de handoff timestamp in its
# "Latest HO Xliff Generate Date" column, so it picks up the same refresh
$overview.Range("G4:G7").Value = "2016-09-07 10:44:28"
